# Auto-generated edit script: refreshes cached market-data values
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) on several
# sheets, as produced by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 570.7375
$ws.Range("J17").Value = 571.3288
$ws.Range("L17").Value = 1713.9864
$ws.Range("N17").Value = -2049.9864
$ws.Range("H33").Value = 1510.6666
$ws.Range("I33").Value = 1030
$ws.Range("J33").Value = 3433.3333
$ws.Range("K33").Value = 1030
$ws.Range("L33").Value = 3433.3333
$ws.Range("M33").Value = -801
$ws.Range("N33").Value = -3891.3333
$ws.Range("H69").Value = 3421.182
$ws.Range("I69").Value = 3479.75
$ws.Range("J69").Value = 3387.7144
$ws.Range("K69").Value = 10439.25
$ws.Range("L69").Value = 10163.1432
$ws.Range("M69").Value = -9565.25
$ws.Range("N69").Value = -11911.1432
$ws.Range("H72").Value = 3421.182
$ws.Range("I72").Value = 3479.75
$ws.Range("J72").Value = 3387.7144
$ws.Range("K72").Value = 31317.75
$ws.Range("L72").Value = 30489.4296
$ws.Range("M72").Value = -26949.75
$ws.Range("N72").Value = -39225.4296
$ws.Range("H76").Value = 3275
$ws.Range("I76").Value = 3281.818
$ws.Range("J76").Value = 3200
$ws.Range("K76").Value = 3281.818
$ws.Range("L76").Value = 3200
$ws.Range("M76").Value = -2966.818
$ws.Range("N76").Value = -3830
$ws.Range("H79").Value = 3275
$ws.Range("I79").Value = 3281.818
$ws.Range("J79").Value = 3200
$ws.Range("K79").Value = 3281.818
$ws.Range("L79").Value = 3200
$ws.Range("M79").Value = -2189.818
$ws.Range("N79").Value = -5384
$ws.Range("H80").Value = 505.85715
$ws.Range("I80").Value = 386.125
$ws.Range("K80").Value = 1158.375
$ws.Range("M80").Value = -160.375
$ws.Range("H83").Value = 505.85715
$ws.Range("I83").Value = 386.125
$ws.Range("K83").Value = 3475.125
$ws.Range("M83").Value = 1516.875
$ws.Range("H88").Value = 61510.8
$ws.Range("I88").Value = 654
$ws.Range("J88").Value = 76725
$ws.Range("K88").Value = 654
$ws.Range("L88").Value = 76725
$ws.Range("M88").Value = -248
$ws.Range("N88").Value = -77537
$ws.Range("H91").Value = 61510.8
$ws.Range("I91").Value = 654
$ws.Range("J91").Value = 76725
$ws.Range("K91").Value = 654
$ws.Range("L91").Value = 76725
$ws.Range("M91").Value = 750
$ws.Range("N91").Value = -79533
$ws.Range("H103").Value = 91759
$ws.Range("I103").Value = 125656.125
$ws.Range("J103").Value = 1366.6666
$ws.Range("K103").Value = 376968.375
$ws.Range("L103").Value = 4099.9998
$ws.Range("M103").Value = -376382.375
$ws.Range("N103").Value = -5271.9998
$ws.Range("H116").Value = 3289.1667
$ws.Range("I116").Value = 2857.1428
$ws.Range("J116").Value = 3894
$ws.Range("K116").Value = 2857.1428
$ws.Range("L116").Value = 3894
$ws.Range("M116").Value = 584.8571999999999
$ws.Range("N116").Value = -10778
$ws.Range("H132").Value = 4237.3145
$ws.Range("I132").Value = 3053.8572
$ws.Range("J132").Value = 8971.143
$ws.Range("K132").Value = 9161.571599999999
$ws.Range("L132").Value = 26913.429
$ws.Range("M132").Value = -6631.571599999999
$ws.Range("N132").Value = -31973.429
$ws.Range("H138").Value = 2472.8333
$ws.Range("I138").Value = 2328.65
$ws.Range("J138").Value = 2575.8215
$ws.Range("K138").Value = 6985.950000000001
$ws.Range("L138").Value = 7727.4645
$ws.Range("M138").Value = -1845.950000000001
$ws.Range("N138").Value = -18007.4645
$ws.Range("H141").Value = 5676.185
$ws.Range("I141").Value = 3081.8333
$ws.Range("J141").Value = 10864.889
$ws.Range("K141").Value = 9245.499899999999
$ws.Range("L141").Value = 32594.667
$ws.Range("M141").Value = -4065.499899999999
$ws.Range("N141").Value = -42954.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 7825.375
$ws.Range("I28").Value = 7825.375
$ws.Range("K28").Value = 7825.375
$ws.Range("M28").Value = -7633.375
$ws.Range("H32").Value = 555530.8
$ws.Range("I32").Value = 664639.9
$ws.Range("J32").Value = 20896.4
$ws.Range("K32").Value = 664639.9
$ws.Range("L32").Value = 20896.4
$ws.Range("M32").Value = -664352.9
$ws.Range("N32").Value = -21470.4
$ws.Range("H61").Value = 3336.0667
$ws.Range("I61").Value = 3038.2856
$ws.Range("J61").Value = 3596.625
$ws.Range("K61").Value = 3038.2856
$ws.Range("L61").Value = 3596.625
$ws.Range("M61").Value = -2826.2856
$ws.Range("N61").Value = -4020.625
$ws.Range("H70").Value = 90000
$ws.Range("J70").Value = 90000
$ws.Range("L70").Value = 90000
$ws.Range("N70").Value = -90540
$ws.Range("H73").Value = 90000
$ws.Range("J73").Value = 90000
$ws.Range("L73").Value = 90000
$ws.Range("N73").Value = -91872
$ws.Range("H74").Value = 2044.0714
$ws.Range("I74").Value = 1337.375
$ws.Range("J74").Value = 2986.3333
$ws.Range("K74").Value = 1337.375
$ws.Range("L74").Value = 2986.3333
$ws.Range("M74").Value = -463.375
$ws.Range("N74").Value = -4734.3333
$ws.Range("H77").Value = 2044.0714
$ws.Range("I77").Value = 1337.375
$ws.Range("J77").Value = 2986.3333
$ws.Range("K77").Value = 6686.875
$ws.Range("L77").Value = 14931.6665
$ws.Range("M77").Value = -2318.875
$ws.Range("N77").Value = -23667.6665
$ws.Range("H99").Value = 7825.375
$ws.Range("I99").Value = 7825.375
$ws.Range("K99").Value = 7825.375
$ws.Range("M99").Value = -4830.375
$ws.Range("H102").Value = 2188.3
$ws.Range("I102").Value = 2098.6667
$ws.Range("J102").Value = 2995
$ws.Range("K102").Value = 2098.6667
$ws.Range("L102").Value = 2995
$ws.Range("M102").Value = -476.6667000000002
$ws.Range("N102").Value = -6239
$ws.Range("H136").Value = 3336.0667
$ws.Range("I136").Value = 3038.2856
$ws.Range("J136").Value = 3596.625
$ws.Range("K136").Value = 9114.856800000001
$ws.Range("L136").Value = 10789.875
$ws.Range("M136").Value = -6564.856800000001
$ws.Range("N136").Value = -15889.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1396.2667
$ws.Range("J80").Value = 179.9
$ws.Range("L80").Value = 179.9
$ws.Range("N80").Value = -2175.9
$ws.Range("H83").Value = 1396.2667
$ws.Range("J83").Value = 179.9
$ws.Range("L83").Value = 899.5
$ws.Range("N83").Value = -10883.5
$ws.Range("H132").Value = 44026.668
$ws.Range("J132").Value = 44026.668
$ws.Range("L132").Value = 44026.668
$ws.Range("N132").Value = -54146.668
$ws.Range("H134").Value = 2803.75
$ws.Range("I134").Value = 2486.238
$ws.Range("J134").Value = 3248.2666
$ws.Range("K134").Value = 7458.714
$ws.Range("L134").Value = 9744.799800000001
$ws.Range("M134").Value = -4923.714
$ws.Range("N134").Value = -14814.7998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3529.1428
$ws.Range("I62").Value = 3176
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3176
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2552
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 3529.1428
$ws.Range("I65").Value = 3176
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 15880
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -12760
$ws.Range("N65").Value = -26240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 754.6667
$ws.Range("I5").Value = 754.6667
$ws.Range("K5").Value = 2264.0001
$ws.Range("M5").Value = -2152.0001
$ws.Range("H23").Value = 45454668
$ws.Range("J23").Value = 66666790
$ws.Range("L23").Value = 200000370
$ws.Range("N23").Value = -200000840
$ws.Range("H113").Value = 1186.15
$ws.Range("J113").Value = 1613.909
$ws.Range("L113").Value = 4841.727000000001
$ws.Range("N113").Value = -9181.727000000001
$ws.Range("H131").Value = 1020.8461
$ws.Range("J131").Value = 1058.9131
$ws.Range("L131").Value = 3176.7393
$ws.Range("N131").Value = -13256.7393
$ws.Range("H134").Value = 5795.839
$ws.Range("I134").Value = 3015.0715
$ws.Range("J134").Value = 8085.8823
$ws.Range("K134").Value = 9045.2145
$ws.Range("L134").Value = 24257.6469
$ws.Range("M134").Value = -3975.2145
$ws.Range("N134").Value = -34397.6469
$ws.Range("H135").Value = 754.6667
$ws.Range("I135").Value = 754.6667
$ws.Range("K135").Value = 6792.0003
$ws.Range("M135").Value = -4257.0003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3089.2246
$ws.Range("I122").Value = 2928.162
$ws.Range("J122").Value = 3585.8333
$ws.Range("K122").Value = 8784.485999999999
$ws.Range("L122").Value = 10757.4999
$ws.Range("M122").Value = -6334.485999999999
$ws.Range("N122").Value = -15657.4999
$ws.Range("H132").Value = 2447.8918
$ws.Range("I132").Value = 1929.7693
$ws.Range("J132").Value = 3672.5454
$ws.Range("K132").Value = 5789.3079
$ws.Range("L132").Value = 11017.6362
$ws.Range("M132").Value = -3259.3079
$ws.Range("N132").Value = -16077.6362
$ws.Range("H136").Value = 3705232.5
$ws.Range("I136").Value = 1664.6666
$ws.Range("J136").Value = 6945854.5
$ws.Range("K136").Value = 4993.9998
$ws.Range("L136").Value = 20837563.5
$ws.Range("M136").Value = -2443.9998
$ws.Range("N136").Value = -20842663.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3473381.5
$ws.Range("I132").Value = 842.12823
$ws.Range("J132").Value = 18521052
$ws.Range("K132").Value = 2526.38469
$ws.Range("L132").Value = 55563156
$ws.Range("M132").Value = 3.615310000000136
$ws.Range("N132").Value = -55568216
$ws.Range("H136").Value = 2020.3273
$ws.Range("I136").Value = 1747.1082
$ws.Range("K136").Value = 5241.3246
$ws.Range("M136").Value = -2691.3246
